$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 (currently "F1"), shifting F1..Z1,Z2 down by one.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new BOM entry.
$ws.Cells.Item(7, 1).Value = "D6,D7"
$ws.Cells.Item(7, 2).Value = "-- mixed values --"
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = "VSSAF522-M3/H"

# Strip the per-cell "left align" style (s="1") from column B and D, leaving
# them with the default (unstyled) cell format - matches the author's edit
# which removed the explicit style index from nearly every B/D cell.
$ws.Range("B1:B24").Style = "Normal"
$ws.Range("D1:D9").Style = "Normal"
$ws.Range("D11:D24").Style = "Normal"

# Update the selection to match the diff.
$ws.Range("H14").Select() | Out-Null
